$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.206.39'
$ws.Range('D3').Value = '1.872.15'
$ws.Range('E3').Value = '  -1.81%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5100'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3766'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07171'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8904'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07597'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.03%  '
$ws.Range('D13').Value = '1.857.83'
$ws.Range('E13').Value = '  -2.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.343'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.45'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008554'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('E18').Value = '  -2.69%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').Value = '27.254.77'
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.082'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = '2.092.47'
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('E24').Value = '  -1.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.13'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.846'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.05'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.127'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.81'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.758'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.724'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08993'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05163'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.093'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7556'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.172'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02043'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.539'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.032'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.078'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5365'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.656'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '114.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.568'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1484'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4680'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.06'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.573'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.73'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.62%  '
